# Update the "想去人数" (F) and "最低票价" (G) values on the 展览 and 全部类型
# sheets to reflect refreshed scrape counts.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1212
$ws1.Range("F4").Value = 297
$ws1.Range("G5").Value = 68
$ws1.Range("F6").Value = 18
$ws1.Range("F7").Value = 12430
$ws1.Range("F8").Value = 64
$ws1.Range("F10").Value = 21
$ws1.Range("F12").Value = 183
$ws1.Range("F13").Value = 12280
$ws1.Range("F14").Value = 4862
$ws1.Range("F15").Value = 4755
$ws1.Range("F16").Value = 147
$ws1.Range("F20").Value = 957
$ws1.Range("F21").Value = 6
$ws1.Range("F23").Value = 172

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1212
$ws4.Range("F4").Value = 297
$ws4.Range("G5").Value = 68
$ws4.Range("F8").Value = 18
$ws4.Range("F9").Value = 12430
$ws4.Range("F10").Value = 64
$ws4.Range("F12").Value = 21
$ws4.Range("F14").Value = 183
$ws4.Range("F15").Value = 12280
$ws4.Range("F16").Value = 4862
$ws4.Range("F17").Value = 4755
$ws4.Range("F18").Value = 147
$ws4.Range("F22").Value = 957
$ws4.Range("F23").Value = 6
$ws4.Range("F25").Value = 172

$wb.Save()
